# Generate Report for Handback
# Populates the "Latest Target File", "Latest Handback File",
# "Latest Handback DateTime" and "Error Detail" columns for the
# 41249a1d-e3a9-478d-979a-85703e472799 row on both the zh-cn and de-de
# sheets, widens the affected columns, and adds a hyperlink on the
# newly populated "Latest Target File" cell.

$wb = $excel.ActiveWorkbook

# Hyperlink-style blue colour (RGB 0x6495ED) used throughout this workbook.
$linkColor = 15570276

# ----------------------------------------------------------------------
# zh-cn sheet
# ----------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Widen columns I (Latest Target File), J (Latest Handback File) and
# P (Error Detail) to 40 characters.
$wsZh.Columns.Item(9).ColumnWidth = 39.17
$wsZh.Columns.Item(10).ColumnWidth = 39.17
$wsZh.Columns.Item(16).ColumnWidth = 39.17

# I5 - Latest Target File: add the handback markdown file name as a hyperlink
# pointing at the current (not-yet-latest) commit of the handback file.
$wsZh.Hyperlinks.Add($wsZh.Range("I5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f231e8996acfc2e22eb62532470c67deec627d32/e2e/41249a1d-e3a9-478d-979a-85703e472799.md", "", "", "41249a1d-e3a9-478d-979a-85703e472799.md")
$wsZh.Range("I5").Font.Name = "Calibri"
$wsZh.Range("I5").Font.Size = 11
$wsZh.Range("I5").Font.Underline = 2
$wsZh.Range("I5").Font.Color = $linkColor

# J5 - Latest Handback File
$wsZh.Range("J5").Value = "41249a1d-e3a9-478d-979a-85703e472799.3ad6a92601e8cf14cada2b5c6064318f0bcae246.zh-cn.xlf"

# K5 - Latest Handback DateTime
$wsZh.Range("K5").Value = "2016-10-20 08:33:59"

# P5 - Error Detail
$wsZh.Range("P5").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f231e8996acfc2e22eb62532470c67deec627d32/e2e/41249a1d-e3a9-478d-979a-85703e472799.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/692d27e8dc467ed550e4188de57e6044d49e7397/e2e/41249a1d-e3a9-478d-979a-85703e472799.md."

# ----------------------------------------------------------------------
# de-de sheet
# ----------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Columns.Item(9).ColumnWidth = 39.17
$wsDe.Columns.Item(10).ColumnWidth = 39.17
$wsDe.Columns.Item(16).ColumnWidth = 39.17

# I5 - Latest Target File
$wsDe.Hyperlinks.Add($wsDe.Range("I5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f231e8996acfc2e22eb62532470c67deec627d32/e2e/41249a1d-e3a9-478d-979a-85703e472799.md", "", "", "41249a1d-e3a9-478d-979a-85703e472799.md")
$wsDe.Range("I5").Font.Name = "Calibri"
$wsDe.Range("I5").Font.Size = 11
$wsDe.Range("I5").Font.Underline = 2
$wsDe.Range("I5").Font.Color = $linkColor

# J5 - Latest Handback File
$wsDe.Range("J5").Value = "41249a1d-e3a9-478d-979a-85703e472799.3ad6a92601e8cf14cada2b5c6064318f0bcae246.de-de.xlf"

# K5 - Latest Handback DateTime
$wsDe.Range("K5").Value = "2016-10-20 08:34:17"

# P5 - Error Detail (same message as zh-cn)
$wsDe.Range("P5").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f231e8996acfc2e22eb62532470c67deec627d32/e2e/41249a1d-e3a9-478d-979a-85703e472799.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/692d27e8dc467ed550e4188de57e6044d49e7397/e2e/41249a1d-e3a9-478d-979a-85703e472799.md."

$wb.Save()
